$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 166, pushing existing rows 166:267 down to 168:269
$ws.Rows("166:167").Insert()

# Fill in the new row 166 (Primera, Región de Arica y Parinacota, $/bandeja 18 kilos)
$ws.Cells.Item(166, 1).Value = 11
$ws.Cells.Item(166, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(166, 3).Value = "Bíobío"
$ws.Cells.Item(166, 4).Value = 44488
$ws.Cells.Item(166, 5).Value = 8
$ws.Cells.Item(166, 6).Value = 100112020
$ws.Cells.Item(166, 7).Value = "Tomate"
$ws.Cells.Item(166, 8).Value = "Larga vida"
$ws.Cells.Item(166, 9).Value = "Primera"
$ws.Cells.Item(166, 10).Value = 600
$ws.Cells.Item(166, 11).Value = 15000
$ws.Cells.Item(166, 12).Value = 16000
$ws.Cells.Item(166, 13).Value = 15500
$ws.Cells.Item(166, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(166, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(166, 16).Value = 861
$ws.Cells.Item(166, 17).Value = 18
$ws.Cells.Item(166, 18).Value = "Hortaliza"

# Fill in the new row 167 (Segunda, Región de Arica y Parinacota, $/bandeja 18 kilos)
$ws.Cells.Item(167, 1).Value = 11
$ws.Cells.Item(167, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(167, 3).Value = "Bíobío"
$ws.Cells.Item(167, 4).Value = 44488
$ws.Cells.Item(167, 5).Value = 8
$ws.Cells.Item(167, 6).Value = 100112020
$ws.Cells.Item(167, 7).Value = "Tomate"
$ws.Cells.Item(167, 8).Value = "Larga vida"
$ws.Cells.Item(167, 9).Value = "Segunda"
$ws.Cells.Item(167, 10).Value = 300
$ws.Cells.Item(167, 11).Value = 14000
$ws.Cells.Item(167, 12).Value = 14000
$ws.Cells.Item(167, 13).Value = 14000
$ws.Cells.Item(167, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(167, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(167, 16).Value = 778
$ws.Cells.Item(167, 17).Value = 18
$ws.Cells.Item(167, 18).Value = "Hortaliza"
